$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency values in column C (rows 2-9)
$ws.Range("C2").Value = 3497
$ws.Range("C3").Value = 3109
$ws.Range("C4").Value = 2127
$ws.Range("C5").Value = 1879
$ws.Range("C6").Value = 1314
$ws.Range("C7").Value = 709
$ws.Range("C8").Value = 613
$ws.Range("C9").Value = 508

# Row 10 changes: category labels swap/change + new frequency
$ws.Range("A10").Value = "Seasonal & Holidays"
$ws.Range("B10").Value = "Home Decor"
$ws.Range("C10").Value = 477

# Row 11 changes: category labels change + new frequency
$ws.Range("A11").Value = "Textiles & Cozy Items"
$ws.Range("B11").Value = "Textiles & Cozy Items"
$ws.Range("C11").Value = 470
